$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'311.08"
$ws.Range("E2").Value = "'0.29%"
$ws.Range("D3").Value = "'37.62"
$ws.Range("E3").Value = "'-1.97%"
$ws.Range("D4").Value = "'5.100"
$ws.Range("E4").Value = "'0.42%"
$ws.Range("D5").Value = "'0.07776"
$ws.Range("E5").Value = "'-2.46%"
$ws.Range("D6").Value = "'4.374"
$ws.Range("E6").Value = "'-1.83%"
$ws.Range("D7").Value = "'1.894"
$ws.Range("E7").Value = "'-6.83%"
$ws.Range("D8").Value = "'8.213"
$ws.Range("E8").Value = "'-1.10%"
$ws.Range("D9").Value = "'2.873"
$ws.Range("E9").Value = "'-7.68%"
$ws.Range("D10").Value = "'0.9221"
$ws.Range("E10").Value = "'-1.17%"
$ws.Range("D11").Value = "'0.1218"
$ws.Range("E11").Value = "'-4.97%"
$ws.Range("D12").Value = "'0.1911"
$ws.Range("E12").Value = "'0.37%"
$ws.Range("D13").Value = "'0.09170"
$ws.Range("E13").Value = "'3.57%"
$ws.Range("D14").Value = "'0.03428"
$ws.Range("E14").Value = "'-0.97%"
$ws.Range("D15").Value = "'0.09683"
$ws.Range("E15").Value = "'-0.11%"
$ws.Range("D16").Value = "'0.001372"
$ws.Range("E16").Value = "'-2.76%"
$ws.Range("D17").Value = "'0.005909"
$ws.Range("E17").Value = "'-7.04%"
$ws.Range("D18").Value = "'3.558"
$ws.Range("E18").Value = "'-0.64%"
$ws.Range("D19").Value = "'0.3403"
$ws.Range("E19").Value = "'0.01%"
$ws.Range("D20").Value = "'5.250"
$ws.Range("E20").Value = "'4.27%"
$ws.Range("E21").Value = "'0.04%"
$ws.Range("D22").Value = "'0.2592"
$ws.Range("E22").Value = "'2.25%"
$ws.Range("D23").Value = "'0.02105"
$ws.Range("E23").Value = "'5,590.74%"
$ws.Range("D24").Value = "'0.04359"
$ws.Range("E24").Value = "'-0.46%"
$ws.Range("E25").Value = "'-2.70%"
$ws.Range("D26").Value = "'0.004256"
$ws.Range("E26").Value = "'-9.02%"
$ws.Range("D27").Value = "'0.0001301"
$ws.Range("E27").Value = "'-63.79%"
$ws.Range("D39").Value = "'0.02093"
$ws.Range("E39").Value = "'-4.14%"
$ws.Range("D40").Value = "'0.05188"
$ws.Range("E40").Value = "'1.74%"
$ws.Range("D41").Value = "'0.007703"
$ws.Range("E41").Value = "'1.21%"
$ws.Range("D42").Value = "'0.009730"
$ws.Range("E42").Value = "'-1.32%"
$ws.Range("D43").Value = "'0.1345"
$ws.Range("E43").Value = "'-2.45%"
$ws.Range("D44").Value = "'0.002062"
$ws.Range("E44").Value = "'0.97%"
$ws.Range("D45").Value = "'0.008908"
$ws.Range("E45").Value = "'0.83%"
$ws.Range("D46").Value = "'0.00006682"
$ws.Range("E46").Value = "'0.36%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'-0.57%"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.001201"
$ws.Range("E48").Value = "'-0.56%"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "'0.002937"
$ws.Range("E49").Value = "'-2.74%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.57%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.57%"
